$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.188.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6016"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07058"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2792"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07625"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.826.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.778"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.000009909"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6254"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.077.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "78.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.831"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.173.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.41%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("E22").Value = "  -5.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.988"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.012"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1297"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.472"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06203"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -15.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.445"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.821"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.791"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.119"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.743"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6394"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.535"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.216.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01735"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.536"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8994"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.990.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.69%  "

$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.583"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05516"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4553"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "

